# Fruta / hortaliza, semanal
# Insert two new weekly records (rows 64 and 65) for "Agrícola del Norte S.A.
# de Arica" - Palta, pushing the existing rows 64-82 down to 66-84.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the two new rows right before the current row 64.
$ws.Rows.Item(64).Insert()
$ws.Rows.Item(64).Insert()

# New row 64: Edranol / Tercera
$ws.Range("A64").Value = 1
$ws.Range("B64").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C64").Value = "Arica y Parinacota"
$ws.Range("D64").Value = 44636
$ws.Range("E64").Value = 15
$ws.Range("F64").Value = "Fruta"
$ws.Range("G64").Value = 100106
$ws.Range("H64").Value = "Oleaginosos"
$ws.Range("I64").Value = 100106002
$ws.Range("J64").Value = "Palta"
$ws.Range("K64").Value = "Edranol"
$ws.Range("L64").Value = "Tercera"
$ws.Range("M64").Value = 250
$ws.Range("N64").Value = 52000
$ws.Range("O64").Value = 55000
$ws.Range("P64").Value = 53500
$ws.Range("Q64").Value = "$/caja 25 kilos"
$ws.Range("R64").Value = "Región de Coquimbo"
$ws.Range("S64").Value = 2140
$ws.Range("T64").Value = 25

# New row 65: Hass / Segunda
$ws.Range("A65").Value = 1
$ws.Range("B65").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C65").Value = "Arica y Parinacota"
$ws.Range("D65").Value = 44636
$ws.Range("E65").Value = 15
$ws.Range("F65").Value = "Fruta"
$ws.Range("G65").Value = 100106
$ws.Range("H65").Value = "Oleaginosos"
$ws.Range("I65").Value = 100106002
$ws.Range("J65").Value = "Palta"
$ws.Range("K65").Value = "Hass"
$ws.Range("L65").Value = "Segunda"
$ws.Range("M65").Value = 200
$ws.Range("N65").Value = 75000
$ws.Range("O65").Value = 78000
$ws.Range("P65").Value = 76500
$ws.Range("Q65").Value = "$/caja 25 kilos"
$ws.Range("R65").Value = "Región de Coquimbo"
$ws.Range("S65").Value = 3060
$ws.Range("T65").Value = 25
